$d = $word.ActiveDocument

# Remove the old _GoBack bookmark (Word will re-add it elsewhere as needed,
# but we handle it explicitly below to match the target location).
foreach ($bm in @($d.Bookmarks)) {
    if ($bm.Name -eq "_GoBack") {
        $bm.Delete()
    }
}

# Replace "three " with "five " (keep trailing space, run-level text will be normalized by Find/Replace)
$d.Content.Find.Execute("three patents, 2 accepted papers, and two research papers under review.", $true, $false, $false, $false, $false, $true, 1, $false, "five patents and published two research papers.", 2)
